# Adds the "Meet Tom" brainstorming paragraphs ahead of the existing
# (bookmark-only) paragraph, and folds the final two new runs into that
# existing paragraph so the _GoBack bookmark stays in place at the end
# of the document content.

$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# The five brand-new paragraphs that land before the original first
# paragraph (built as whole <w:p> fragments so Word's XML importer keeps
# each <w:r> distinct instead of merging same-format runs together).
$newParasXml = `
  ('<w:p ' + $wns + '>' + `
     '<w:r><w:t xml:space="preserve">Something I want to improve is how to conduct an interview. </w:t></w:r>' + `
     '<w:r><w:t xml:space="preserve">There are so many awesome websites out there to look for the </w:t></w:r>' + `
     '<w:r><w:t>information, but the material is not in</w:t></w:r>' + `
   '</w:p>') + `
  ('<w:p ' + $wns + '/>') + `
  ('<w:p ' + $wns + '/>') + `
  ('<w:p ' + $wns + '>' + `
     '<w:r><w:t>As a user I can put skills I want to learn so that people can teach or show me the correct path to take in order to achieve that skills.</w:t></w:r>' + `
   '</w:p>') + `
  ('<w:p ' + $wns + '>' + `
     '<w:r><w:t xml:space="preserve">As a user I can put the skills I want to </w:t></w:r>' + `
   '</w:p>')

$start = $d.Range(0, 0)
$start.InsertXML($newParasXml)

# Locate the original paragraph again via its still-intact _GoBack
# bookmark, then rewrite that paragraph so the two new runs sit in front
# of the (re-emitted) bookmark markers.
$bm = $d.Bookmarks.Item("_GoBack")
$targetPara = $bm.Range.Paragraphs.Item(1)
$targetRange = $targetPara.Range

$lastParaXml = ('<w:p ' + $wns + '>' + `
     '<w:r><w:t xml:space="preserve">I signed up for goalbook </w:t></w:r>' + `
     '<w:r><w:t xml:space="preserve">. My goal is to </w:t></w:r>' + `
     '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
     '<w:bookmarkEnd w:id="0"/>' + `
   '</w:p>')

$targetRange.InsertXML($lastParaXml)
